$d = $word.ActiveDocument

# 1. Update the "Curso (semestre ideal)" line: drop the "EQD (10), " part.
$d.Content.Find.Execute(
    "Curso (semestre ideal): EQD (10), EQN (12)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Curso (semestre ideal): EQN (12)", 2
)

# 2. Append a new "Requisitos" section at the end of the document:
#      Heading2  "Requisitos"
#      ListBullet "LOQ4044 -  Introdução à Engenharia da Qualidade  (Requisito fraco)" + line break

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$headingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$headingPara.Style = "Heading 2"
$headingPara.Range.Text = "Requisitos"

$headingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$headingPara.Range.InsertParagraphAfter()

$bulletPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bulletPara.Style = "List Bullet"
$bulletPara.Range.Text = "LOQ4044 -  Introdução à Engenharia da Qualidade  (Requisito fraco)"
$bulletPara.Range.InsertAfter([char]11)
